# Updated cryptos list (GitHub Actions bot refresh of Price / Volume(1h) columns).
#
# Notes:
#  - Price (column D) cells hold plain text in the source sheet (e.g. "27.508.19",
#    "1.006"), not numbers. Excel auto-detects single-dot, all-digit strings
#    ("1.004", "226.13", ...) as numeric when assigned via .Value, which would
#    both change the cell's stored type and silently drop significant trailing
#    zeros (e.g. "0.07730" -> 0.0773). To keep these as literal text -- matching
#    the source data -- a leading apostrophe is used for every new Price value
#    that looks numeric, the same trick used when typing such values directly
#    into Excel.
#  - Rows 41/42 swapped coin order (Maker now ranks above TrustWalletToken), so
#    their Coin/Link/Price/Volume cells are fully rewritten rather than just
#    the Price/Volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.522.35"
$ws.Range("E2").Value = "  +5.67%  "

$ws.Range("D3").Value = "1.723.15"

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'226.13"
$ws.Range("E5").Value = "  +3.63%  "

$ws.Range("D6").Value = "'0.5341"
$ws.Range("E6").Value = "  +2.86%  "

$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "'0.2669"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("D9").Value = "'0.06591"
$ws.Range("E9").Value = "  +4.30%  "

$ws.Range("D10").Value = "'21.69"
$ws.Range("E10").Value = "  +6.48%  "

$ws.Range("D11").Value = "'0.07730"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").Value = "1.720.28"
$ws.Range("E13").Value = "  +4.24%  "

$ws.Range("D14").Value = "1.960.09"
$ws.Range("E14").Value = "  +4.46%  "

$ws.Range("D15").Value = "'0.5825"
$ws.Range("E15").Value = "  +4.28%  "

$ws.Range("D16").Value = "0.0₅8300"
$ws.Range("E16").Value = "  +2.12%  "

$ws.Range("D17").Value = "'67.90"
$ws.Range("E17").Value = "  +4.16%  "

$ws.Range("D18").Value = "27.521.29"

$ws.Range("D19").Value = "'219.99"
$ws.Range("E19").Value = "  +14.94%  "

$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").Value = "'4.734"
$ws.Range("E21").Value = "  +2.64%  "

$ws.Range("D22").Value = "'10.65"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").Value = "'6.086"
$ws.Range("E23").Value = "  +3.13%  "

$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").Value = "'147.40"
$ws.Range("E25").Value = "  +2.74%  "

$ws.Range("D26").Value = "'1.738"
$ws.Range("E26").Value = "  +15.27%  "

$ws.Range("D27").Value = "'0.1235"
$ws.Range("E27").Value = "  +4.40%  "

$ws.Range("D28").Value = "'7.408"
$ws.Range("E28").Value = "  +3.13%  "

$ws.Range("D29").Value = "'16.56"
$ws.Range("E29").Value = "  +4.47%  "

$ws.Range("D30").Value = "'0.05556"
$ws.Range("E30").Value = "  +3.47%  "

$ws.Range("E31").Value = "  +2.78%  "

$ws.Range("D32").Value = "'3.565"
$ws.Range("E32").Value = "  +3.27%  "

$ws.Range("D33").Value = "'3.446"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("D34").Value = "'1.665"
$ws.Range("E34").Value = "  +7.28%  "

$ws.Range("E35").Value = "  +3.03%  "

$ws.Range("D36").Value = "'0.9653"
$ws.Range("E36").Value = "  +2.12%  "

$ws.Range("D37").Value = "'2.422"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'0.5961"
$ws.Range("E38").Value = "  +5.93%  "

$ws.Range("D39").Value = "'0.01655"
$ws.Range("E39").Value = "  +4.91%  "

$ws.Range("D40").Value = "'5.915"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.056.24"
$ws.Range("E41").Value = "  +2.57%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8534"
$ws.Range("E42").Value = "  +3.52%  "

$ws.Range("D43").Value = "'1.005"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").Value = "'101.60"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").Value = "1.867.09"
$ws.Range("E45").Value = "  +4.53%  "

$ws.Range("E46").Value = "  +3.92%  "

$ws.Range("D47").Value = "'58.98"
$ws.Range("E47").Value = "  +3.08%  "

$ws.Range("D48").Value = "'8.223"
$ws.Range("E48").Value = "  +3.96%  "

$ws.Range("D49").Value = "'0.4439"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").Value = "'0.05244"
$ws.Range("E51").Value = "  +2.02%  "
